# "Add new periods in dict"
# Adds 10 new period rows (Jul/Aug/Sep 2024 for the MAT / 2MATs / 3MMT / Month
# period types) to the period_lbl sheet, highlights them in yellow, fixes up
# the sheet's used range / column widths, re-applies the sort that produced
# the final row order, and switches the active sheet/tab back to period_lbl.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("period_lbl")

# New rows to append right after the existing data (row 111 = A:110).
# Tuple layout: row, A (label_num), B (period_lbl), C (period_code)
$newRows = @(
    @(112, 111, "2MATs: Sep 2024",  "2MATs: 2024 (09) Sep"),
    @(113, 112, "3MMT: Aug 2024",   "3MMT: 2024 (07) Jul"),
    @(114, 113, "3MMT: Jul 2024",   "3MMT: 2024 (08) Aug"),
    @(115, 114, "3MMT: Sep 2024",   "3MMT: 2024 (09) Sep"),
    @(116, 115, "MAT: Jul 2024",    "MAT: 2024 (07) Jul"),
    @(117, 116, "MAT: Aug 2024",    "MAT: 2024 (08) Aug"),
    @(118, 117, "MAT: Sep 2024",    "MAT: 2024 (09) Sep"),
    @(119, 118, "Month: Jul 2024",  "Month: 2024 (07) Jul"),
    @(120, 119, "Month: Aug 2024",  "Month: 2024 (08) Aug"),
    @(121, 120, "Month: Sep 2024",  "Month: 2024 (09) Sep")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]

    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3)).Interior.Color = 65535
}

# Row 118 ("MAT: Sep 2024") keeps its label text quote-prefixed, same as the
# author's original entry.
$ws.Cells.Item(118, 2).Value = "'MAT: Sep 2024"
$ws.Cells.Item(118, 2).Interior.Color = 65535

# Re-create the sort the author ran over the new rows (column A then D) so the
# worksheet keeps the same <sortState> marker.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A113:A121")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("D113:D121")) | Out-Null
$ws.Sort.SetRange($ws.Range("A113:D121"))
$ws.Sort.Apply()

# Column width tweaks that came with the new data.
$ws.Columns.Item(3).ColumnWidth = 12.25
$ws.Columns.Item(7).ColumnWidth = 29.45

# Make period_lbl the active/selected sheet & tab again (it had been on
# time_period_type before this edit).
$ws.Activate()
$ws.Range("G115").Select()
